$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp (shared string used by A1) ---
$ws.Range("A1").Value = "Datos actualizados a 29 de Marzo de 2020 a las 20:20"

# --- Estados Unidos (row 4) gets refreshed totals ---
$ws.Range("B4").Value = 135510
$ws.Range("C4").Value = 11932
$ws.Range("D4").Value = 4378
$ws.Range("E4").Value = 128748
$ws.Range("F4").Value = 2948
$ws.Range("G4").Value = 164
$ws.Range("H4").Value = 2384

# --- Japon is re-sorted up, right after Ecuador (was after Rumania) ---
# Row 32 becomes Japon (fresh numbers), Polonia/Rumania each shift down one row
# keeping their previous values.
$ws.Range("A32").Value = "Japon"
$ws.Range("B32").Value = 1866
$ws.Range("C32").Value = 173
$ws.Range("D32").Value = 424
$ws.Range("E32").Value = 1388
$ws.Range("F32").Value = 56
$ws.Range("G32").Value = 2
$ws.Range("H32").Value = 54

$ws.Range("A33").Value = "Polonia"
$ws.Range("B33").Value = 1771
$ws.Range("C33").Value = 133
$ws.Range("D33").Value = 7
$ws.Range("E33").Value = 1744
$ws.Range("F33").Value = 3
$ws.Range("G33").Value = 2
$ws.Range("H33").Value = 20

$ws.Range("A34").Value = "Rumania"
$ws.Range("B34").Value = 1760
$ws.Range("C34").Value = 308
$ws.Range("D34").Value = 169
$ws.Range("E34").Value = 1549
$ws.Range("F34").Value = 34
$ws.Range("G34").Value = 5
$ws.Range("H34").Value = 42

# --- Burkina Faso is re-sorted up, right after San Marino (was after Azerbaiyan) ---
# Row 87 becomes Burkina Faso (fresh numbers), Republica de Chipre/Albania/Azerbaiyan
# each shift down one row keeping their previous values.
$ws.Range("A87").Value = "Burkina Faso"
$ws.Range("B87").Value = 222
$ws.Range("C87").Value = 15
$ws.Range("D87").Value = 23
$ws.Range("E87").Value = 187
$ws.Range("F87").Value = 0
$ws.Range("G87").Value = 1
$ws.Range("H87").Value = 12

$ws.Range("A88").Value = "Republica de Chipre"
$ws.Range("B88").Value = 214
$ws.Range("C88").Value = 35
$ws.Range("D88").Value = 15
$ws.Range("E88").Value = 194
$ws.Range("F88").Value = 3
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 5

$ws.Range("A89").Value = "Albania"
$ws.Range("B89").Value = 212
$ws.Range("C89").Value = 15
$ws.Range("D89").Value = 33
$ws.Range("E89").Value = 169
$ws.Range("F89").Value = 3
$ws.Range("G89").Value = 0
$ws.Range("H89").Value = 10

$ws.Range("A90").Value = "Azerbaiyan"
$ws.Range("B90").Value = 209
$ws.Range("C90").Value = 27
$ws.Range("D90").Value = 15
$ws.Range("E90").Value = 190
$ws.Range("F90").Value = 23
$ws.Range("G90").Value = 0
$ws.Range("H90").Value = 4
